$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update summary header values
# ---------------------------------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 5666116
# Cant. Trabajadores
$ws.Range("C13").Value = 6
# Cant. Periodos
$ws.Range("F13").Value = 45

# ---------------------------------------------------------------------------
# 2) Grow the detail table from 82 data rows (16-97) to 88 data rows (16-103).
#    Insert 6 new rows right after the current last data row (97) and before
#    the blank spacer rows, so the trailing signature block shifts from
#    rows 102-103 down to rows 108-109.
# ---------------------------------------------------------------------------
$ws.Rows("98:103").Insert()

# Preserve the special "last row" bottom-border style (currently still on
# row 97) by stashing a copy of it far away before re-styling everything.
$ws.Range("B97:J97").Copy()
$ws.Range("B500:J500").PasteSpecial(-4122)

# Re-apply the normal interior-row style across the whole (now larger) table
$ws.Range("B16:J16").Copy()
$ws.Range("B16:J103").PasteSpecial(-4122)

# Put the special last-row style back onto the new final row (103)
$ws.Range("B500:J500").Copy()
$ws.Range("B103:J103").PasteSpecial(-4122)

# Clean up the scratch cells used to stash the style
$ws.Range("B500:J500").Clear()

# ---------------------------------------------------------------------------
# 3) Rewrite the detail rows (B16:G103) with the refreshed dataset: two new
#    workers (EDSON GUZMAN ALTAMAR, MARLON ANTONIO LONG MATOS) and one more
#    period (2508) added for every active worker.
# ---------------------------------------------------------------------------
$tableData = @(
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2112",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2201",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2202",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2203",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2204",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2205",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2206",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2207",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2208",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2209",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2210",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2211",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2212",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2301",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2302",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2303",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2304",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2305",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2306",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2307",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2308",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2309",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2310",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2311",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2312",40972,1024300)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2401",40972,1024300)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2402",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2402",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2402",42224,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2403",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2403",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2403",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2404",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2404",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2404",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2405",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2405",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2405",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2406",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2406",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2406",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2407",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2407",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2407",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2408",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2408",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2408",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2409",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2409",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2409",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2410",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2410",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2410",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2411",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2411",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2411",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2412",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2412",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2412",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2501",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2501",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2501",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2502",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2502",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2502",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2503",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2503",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2503",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2504",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2504",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2504",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2505",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2505",40972,1024300)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2505",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2506",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2506",40972,1024300)
    ,("CC","1010095635","STEPHANY PAOLA BARRIOS EMILIANI","2506",51440,1929000)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2506",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2507",84120,2103000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2507",40972,1024300)
    ,("CC","1010095635","STEPHANY PAOLA BARRIOS EMILIANI","2507",77160,1929000)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2507",105560,2639000)
    ,("CC","1047454218","FABER ANTONIO MURIEL VALDES","2508",84120,2103000)
    ,("CC","1007576861","EDSON GUZMAN ALTAMAR","2508",21120,1584000)
    ,("CC","1047442262","MISAEL DAVID CORPAS MARTINEZ","2508",40972,1024300)
    ,("CC","1143365346","MARLON ANTONIO LONG MATOS","2508",54912,1584000)
    ,("CC","1010095635","STEPHANY PAOLA BARRIOS EMILIANI","2508",77160,1929000)
    ,("CC","1033717238","LADY JOHANA SILVA CASTILLO","2508",105560,2639000)
)

$startRow = 16
$r = $startRow
foreach ($row in $tableData) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
